$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.073.60"
$ws.Range("E2").Value = "  +3.33%  "
$ws.Range("D3").Value = "1.690.13"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'216.72"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'24.15"
$ws.Range("E8").Value = "  +5.68%  "
$ws.Range("E9").Value = "  +1.94%  "
$ws.Range("D10").Value = "'0.0626"
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("D11").Value = "'0.0885"
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").Value = "1.930.50"
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("D13").Value = "1.689.79"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").Value = "'0.558"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").Value = "'66.91"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "'250.58"
$ws.Range("E17").Value = "  +6.30%  "
$ws.Range("D18").Value = "28.034.90"
$ws.Range("E18").Value = "  +3.20%  "
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").Value = "'7.69"
$ws.Range("E20").Value = "  -3.40%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "'4.54"
$ws.Range("E22").Value = "  -0.37%  "
$ws.Range("D23").Value = "'9.55"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "'2.05"
$ws.Range("E24").Value = "  -1.93%  "
$ws.Range("D25").Value = "'147.73"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.24"
$ws.Range("E30").Value = "  +5.94%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.0504"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("E33").Value = "  -1.86%  "
$ws.Range("D34").Value = "1.446.67"
$ws.Range("E34").Value = "  -6.61%  "
$ws.Range("E35").Value = "  -2.94%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("E38").Value = "  -1.80%  "
$ws.Range("D39").Value = "'0.0173"
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("E40").Value = "  -2.04%  "
$ws.Range("D41").Value = "'69.54"
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").Value = "'5.55"
$ws.Range("E43").Value = "  -3.43%  "
$ws.Range("D44").Value = "1.837.75"
$ws.Range("E44").Value = "  +0.53%  "
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("D46").Value = "'0.799"
$ws.Range("E46").Value = "  +1.08%  "
$ws.Range("E47").Value = "  +7.36%  "
$ws.Range("D48").Value = "'89.57"
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("E49").Value = "  -1.30%  "
$ws.Range("E50").Value = "  -0.81%  "
$ws.Range("D51").Value = "'7.94"
$ws.Range("E51").Value = "  -3.79%  "

# Reset style on cells that Excel auto-marked with quotePrefix due to numeric-looking text,
# so they match the original plain (unstyled) inlineStr cells.
$ws.Range("D5").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").Style = "Normal"
